$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the question text in C3 (keep insertion-order so the shared-string
#     table comes out in the same order as the authored workbook) ---
$ws.Range("C3").Value = "Where is Anodiam  located (Address & Phone Number) ?"

# --- Add new rows 16-18 ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Parent"
$ws.Range("C16").Value = "Is Anodiam providing labs for AI ?"
$ws.Range("D16").Value = "Y"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Parent"
$ws.Range("C17").Value = "Is there any internal tests ?"
$ws.Range("D17").Value = "Y"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Parent"
$ws.Range("C18").Value = "Do Anodiam have any App ?"
$ws.Range("D18").Value = "Y"

# --- Change the header label for column B from "Student" to "User" ---
$ws.Range("B1").Value = "User"

# --- Add new rows 19-23 ---
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Parent"
$ws.Range("C19").Value = "Do Anodiam provide course materials ?"
$ws.Range("D19").Value = "Y"

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Parent"
$ws.Range("C20").Value = "Do Anodiam conduct online classes ?"
$ws.Range("D20").Value = "Y"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Parent"
$ws.Range("C21").Value = "After taking the course what is the prospect of my son/daughter ?"
$ws.Range("D21").Value = "N"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "Principal"
$ws.Range("C22").Value = "Will Anodiam conduct classes in school premises ?"
$ws.Range("D22").Value = "Y"

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "Principal"
$ws.Range("C23").Value = "What are the extra things Anodiam will provide to students other than AI ?"
$ws.Range("D23").Value = "Y"

# --- Widen column B slightly to fit the new "Principal" values ---
$ws.Columns.Item(2).ColumnWidth = 6.83

# --- Replace the data validation lists so they cover the whole column and
#     include the new roles / Y-N dropdown ---
$bValidation = $ws.Range("B1:B1048576")
$bValidation.Validation.Delete()
$bValidation.Validation.Add(3, 1, 1, """Parent,Student,Principal,Teacher""")

$dValidation = $ws.Range("D1:D1048576")
$dValidation.Validation.Add(3, 1, 1, """Y,N""")

# --- Move the active selection to reflect where the author left off editing ---
[void]$ws.Range("A24:A25").Select()
